$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 ("business card")
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Shape 2: title textbox "Create a DC cycling safety map collectively in
# just 2 steps" -> "Help us create DC's cycling safety map in just 2 steps"
$titleShape = $s1.Shapes.Item(2)
$apostrophe = [char]0x2019
$titleShape.TextFrame.TextRange.Text = "Help us create DC" + $apostrophe + "s cycling safety map in just 2 steps"
$titleShape.Width = 26792847 / 12700

# Shape 5: "Rate whenever you want, even with your mobile device" ->
# "We will raffle twelve $25 gift cards among all participants. "
$raffleShape = $s1.Shapes.Item(5)
$raffleRange = $raffleShape.TextFrame.TextRange
$raffleRange.Delete()
[void]$raffleRange.InsertAfter('We will raffle twelve $25 gift cards among all participants. ')
$raffleShape.Left = 1540215 / 12700
$raffleShape.Width = 29323733 / 12700

# Shape 6: "Starts Now at https://www.cyclingsafety.umd.edu" ->
# "Participate: https://www.cyclingsafety.umd.edu"
$participateShape = $s1.Shapes.Item(6)
$participateRange = $participateShape.TextFrame.TextRange
$participateRange.Characters(1, 10).Text = "Participate:"
$participateRange.Characters(13, 1).Delete()
$participateRange.Characters(13, 3).Text = " "
$participateShape.Top = 11513190 / 12700
$participateShape.Width = 25762950 / 12700

# Shape 8: UMD Urban Computing Lab logo picture, shifted up
$logoPic = $s1.Shapes.Item(8)
$logoPic.Top = 11737829 / 12700

# ---------------------------------------------------------------------------
# Slide 2 ("Make our city more bike-friendly")
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$bikeShape = $s2.Shapes.Item(1)
$bikeShape.TextFrame.WordWrap = -1
$bikeShape.TextFrame.TextRange.Font.Size = 130
$bikeShape.TextFrame.TextRange.Text = "Help make our cities more bike-friendly"
$bikeShape.Left = 2228349 / 12700
$bikeShape.Top = 8059087 / 12700
$bikeShape.Width = 27547301 / 12700
$bikeShape.Height = 2169825 / 12700
